$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 70/71. Every existing row from 70 downward
# (70..172) shifts down by two (to 72..174), which is exactly the shift
# described by the diff (new row N = old row N-2 for N >= 72, and the
# two brand-new records land at the freshly opened 70/71).
$ws.Rows("70:71").Insert()

# New record #1 ("Primera" quality) -> row 70
$ws.Range("A70").Value = 1
$ws.Range("B70").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C70").Value = "Arica y Parinacota"
$ws.Range("D70").Value = 44477
$ws.Range("E70").Value = 15
$ws.Range("F70").Value = 100114014
$ws.Range("G70").Value = "Betarraga"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 1200
$ws.Range("K70").Value = 450
$ws.Range("L70").Value = 500
$ws.Range("M70").Value = 475
$ws.Range("N70").Value = "$/paquete 4 unidades"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 119
$ws.Range("Q70").Value = 4
$ws.Range("R70").Value = "Hortaliza"

# New record #1 ("Segunda" quality) -> row 71
$ws.Range("A71").Value = 1
$ws.Range("B71").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C71").Value = "Arica y Parinacota"
$ws.Range("D71").Value = 44477
$ws.Range("E71").Value = 15
$ws.Range("F71").Value = 100114014
$ws.Range("G71").Value = "Betarraga"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Segunda"
$ws.Range("J71").Value = 1400
$ws.Range("K71").Value = 450
$ws.Range("L71").Value = 500
$ws.Range("M71").Value = 475
$ws.Range("N71").Value = "$/paquete 5 unidades"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 95
$ws.Range("Q71").Value = 5
$ws.Range("R71").Value = "Hortaliza"
